$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/date range) ---
$ws.Range("A8").Value = "Volume 30   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/30/2023  Through  11/5/2023"

# --- Simple numeric value updates (style/type unchanged) ---
$ws.Range("N14").Value = -79.487179487179
$ws.Range("G15").Value = 7
$ws.Range("J15").Value = 28
$ws.Range("K15").Value = -28.571428571428
$ws.Range("L15").Value = -25.925925925925
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -18.181818181818
$ws.Range("F16").Value = 28
$ws.Range("G16").Value = 26
$ws.Range("H16").Value = 7.692307692307
$ws.Range("I16").Value = 284
$ws.Range("J16").Value = 260
$ws.Range("K16").Value = 9.230769230769
$ws.Range("L16").Value = 82.051282051282
$ws.Range("M16").Value = 17.842323651452
$ws.Range("N16").Value = -61.042524005487
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = -28.571428571428
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 39
$ws.Range("H17").Value = -5.128205128205
$ws.Range("I17").Value = 414
$ws.Range("J17").Value = 369
$ws.Range("K17").Value = 12.195121951219
$ws.Range("L17").Value = 56.226415094339
$ws.Range("M17").Value = 24.698795180722
$ws.Range("N17").Value = -6.966292134831
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 400
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 22.222222222222
$ws.Range("I18").Value = 148
$ws.Range("J18").Value = 134
$ws.Range("K18").Value = 10.447761194029
$ws.Range("L18").Value = 78.313253012048
$ws.Range("M18").Value = -21.276595744680
$ws.Range("N18").Value = -83.718371837183
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 200
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = -8.823529411764
$ws.Range("I19").Value = 388
$ws.Range("J19").Value = 347
$ws.Range("K19").Value = 11.815561959654
$ws.Range("L19").Value = 39.568345323741
$ws.Range("M19").Value = 70.175438596491
$ws.Range("N19").Value = 0
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = -16.666666666666
$ws.Range("I20").Value = 262
$ws.Range("J20").Value = 179
$ws.Range("K20").Value = 46.368715083798
$ws.Range("L20").Value = 123.931623931624
$ws.Range("M20").Value = 291.044776119403
$ws.Range("N20").Value = -27.423822714681
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 122
$ws.Range("G21").Value = 134
$ws.Range("H21").Value = -8.955223880597
$ws.Range("I21").Value = 1524
$ws.Range("J21").Value = 1322
$ws.Range("K21").Value = 15.279878971255
$ws.Range("L21").Value = 62.820512820512
$ws.Range("M21").Value = 41.111111111111
$ws.Range("N21").Value = -47.430148327009
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -25
$ws.Range("J23").Value = 16
$ws.Range("K23").Value = 68.75
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = 62.5
$ws.Range("F24").Value = 78
$ws.Range("G24").Value = 46
$ws.Range("H24").Value = 69.565217391304
$ws.Range("I24").Value = 777
$ws.Range("J24").Value = 670
$ws.Range("K24").Value = 15.970149253731
$ws.Range("L24").Value = 73.051224944320
$ws.Range("M24").Value = 42.568807339449
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 60
$ws.Range("F25").Value = 62
$ws.Range("H25").Value = 44.186046511627
$ws.Range("I25").Value = 490
$ws.Range("J25").Value = 462
$ws.Range("K25").Value = 6.060606060606
$ws.Range("L25").Value = 47.147147147147
$ws.Range("M25").Value = -14.035087719298
$ws.Range("D26").Value = 1
$ws.Range("G26").Value = 10
$ws.Range("H26").Value = -90
$ws.Range("J26").Value = 49
$ws.Range("K26").Value = -18.367346938775
$ws.Range("L26").Value = -27.272727272727
$ws.Range("C27").Value = 4
$ws.Range("F27").Value = 13
$ws.Range("H27").Value = 550
$ws.Range("I27").Value = 88
$ws.Range("K27").Value = 10
$ws.Range("L27").Value = -4.347826086956
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 21
$ws.Range("K28").Value = -8.695652173913
$ws.Range("L28").Value = -36.363636363636
$ws.Range("M28").Value = -12.5
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 19
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = -24
$ws.Range("M29").Value = -9.523809523809
$ws.Range("N29").Value = -77.906976744186

# --- Cells switching to TEXT type with style 14 (shared "0" / "***.*" placeholders) ---
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("F14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("F14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("F14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

# --- Cells switching to NUMERIC type with style 15 ---
$ws.Range("G14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1
$ws.Range("G14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 1

$excel.CutCopyMode = $false
